# Update column F (dSF) values on Sheet1 to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -1
    3  = -2
    4  = -6
    5  = 2
    6  = -2
    8  = 2
    10 = -2
    11 = -2
    12 = 4
    14 = 1
    15 = 4
    16 = -4
    17 = 5
    18 = 1
    19 = 2
    21 = -2
    22 = 2
    24 = -4
    25 = -6
    26 = 1
    27 = 7
    28 = -2
    29 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
